$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset formatting across the data block so re-written rows start clean
$ws.Range("A3:D56").ClearFormats()

# --- Country mapping table rows (COSIMO ISO3/Country, GENUS ISO3/Country) ---
$rows = @(
    @("AFG", "Afghanistan", "PAK", "Pakistan"),
    @("ANT", "Netherlands Antilles", "VEN", "Venezuela"),
    @("BDI", "Burundi", "RWA", "Rwanda"),
    @("BHR", "Bahrain", "SAU", "Saudi Arabia"),
    @("BLX", "Belgium-Luxembourg", "BEL", "Belgium"),
    @("BMU", "Bermuda", "USA", "United States"),
    @("BTN", "Bhutan", "NPL", "Nepal"),
    @("COD", "Congo - Kinshasa", "COG", "Congo"),
    @("COM", "Comoros", "MDG", "Madagascar"),
    @("CZ2", "Czechoslovakia", "CZE", "Czech Republic"),
    @("DMA", "Dominica", "LCA", "Saint Lucia"),
    @("ERI", "Eritrea", "DJI", "Djibouti"),
    @("ESH", "Western Sahara", "MAR", "Morocco"),
    @("ET2", "Ethiopia PDR", "ETH", "Ethiopia"),
    @("FSM", "Micronesia (Federated States of)", "FJI", "Fiji"),
    @("GAB", "Gabon", "COG", "Congo"),
    @("GNQ", "Equatorial Guinea", "CMR", "Cameroon"),
    @("HKG", "Hong Kong SAR China", "CHN", "China"),
    @("KHM", "Cambodia", "THA", "Thailand"),
    @("KIR", "Kiribati", "PYF", "French Polynesia"),
    @("KNA", "St. Kitts & Nevis", "ATG", "Antigua and Barbuda"),
    @("LBR", "Liberia", "CIV", "Côte d’Ivoire"),
    @("LSO", "Lesotho", "ZAF", "South Africa"),
    @("MAC", "Macau SAR China", "CHN", "China"),
    @("MHL", "Marshall Islands", "FJI", "Fiji"),
    @("MMR", "Myanmar (Burma)", "LAO", "Laos"),
    @("OMN", "Oman", "YEM", "Yemen"),
    @("PLW", "Palau", "PHL", "Philippines"),
    @("PNG", "Papua New Guinea", "IDN", "Indonesia"),
    @("PRI", "Puerto Rico", "CUB", "Cuba"),
    @("PRK", "North Korea", "KOR", "South Korea"),
    @("QAT", "Qatar", "ARE", "United Arab Emirates"),
    @("SGP", "Singapore", "MYS", "Malaysia"),
    @("SLB", "Solomon Islands", "NCL", "New Caledonia"),
    @("SLE", "Sierra Leone", "GIN", "Guinea"),
    @("SMR", "San Marino", "ITA", "Italy"),
    @("SOM", "Somalia", "ETH", "Ethiopia"),
    @("SRM", "Serbia and Montenegro", "SRB", "Serbia"),
    @("STP", "São Tomé and Príncipe", "CMR", "Cameroon"),
    @("SYC", "Seychelles", "MDG", "Madagascar"),
    @("TCD", "Chad", "SDN", "Sudan"),
    @("TGO", "Togo", "BEN", "Benin"),
    @("TKM", "Turkmenistan", "UZB", "Uzbekistan"),
    @("TLS", "Timor-Leste", "IDN", "Indonesia"),
    @("TON", "Tonga", "FJI", "Fiji"),
    @("TUV", "Tuvalu", "FJI", "Fiji"),
    @("TWN", "Taiwan", "CHN", "China"),
    @("UGA", "Uganda", "KEN", "Kenya"),
    @("USR", "USSR", "RUS", "Russia"),
    @("VNM", "Vietnam", "LAO", "Laos"),
    @("VUT", "Vanuatu", "FJI", "Fiji"),
    @("WSM", "Samoa", "FJI", "Fiji"),
    @("YUG", "Yugoslav SFR", "SRB", "Serbia"),
    @("ZMB", "Zambia", "ZWE", "Zimbabwe")
)

$startRow = 3
for ($i = 0; $i -lt $rows.Length; $i++) {
    $row = $rows[$i]
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# --- Re-apply the "pasted" font (Lucida Grande 11, black) on the specific cells ---
$pastedCells = @("B9","B15","B19","D24","D37","D38","D45")
foreach ($addr in $pastedCells) {
    $ws.Range($addr).Font.Name = "Lucida Grande"
    $ws.Range($addr).Font.Size = 11
    $ws.Range($addr).Font.Color = 0
}

# --- Column B widened to fit the longer country names now present ---
$ws.Columns("B").ColumnWidth = 28.6640625

# --- Restore selection to match the last-edited cell ---
$ws.Range("C39").Select()
